$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 321.30768
$ws.Range("I2").Value = 323.16666
$ws.Range("K2").Value = 323.16666
$ws.Range("M2").Value = -210.16666
$ws.Range("H40").Value = 1888.1923
$ws.Range("I40").Value = 1643.8889
$ws.Range("J40").Value = 2017.5294
$ws.Range("K40").Value = 1643.8889
$ws.Range("L40").Value = 2017.5294
$ws.Range("M40").Value = -1468.8889
$ws.Range("N40").Value = -2367.5294
$ws.Range("H41").Value = 1395.5385
$ws.Range("I41").Value = 1695.1
$ws.Range("K41").Value = 1695.1
$ws.Range("M41").Value = -1255.1
$ws.Range("H69").Value = 7472
$ws.Range("J69").Value = 9259
$ws.Range("L69").Value = 27777
$ws.Range("N69").Value = -29525
$ws.Range("H72").Value = 7472
$ws.Range("J72").Value = 9259
$ws.Range("L72").Value = 83331
$ws.Range("N72").Value = -92067
$ws.Range("H86").Value = 7360.364
$ws.Range("I86").Value = 7432.5
$ws.Range("J86").Value = 7344.3335
$ws.Range("K86").Value = 7432.5
$ws.Range("L86").Value = 7344.3335
$ws.Range("M86").Value = -6309.5
$ws.Range("N86").Value = -9590.333500000001
$ws.Range("H88").Value = 6500
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 6500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 6500
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -7312
$ws.Range("H89").Value = 7360.364
$ws.Range("I89").Value = 7432.5
$ws.Range("J89").Value = 7344.3335
$ws.Range("K89").Value = 37162.5
$ws.Range("L89").Value = 36721.6675
$ws.Range("M89").Value = -31546.5
$ws.Range("N89").Value = -47953.6675
$ws.Range("H91").Value = 6500
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 6500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 6500
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -9308
$ws.Range("H106").Value = 1198.4667
$ws.Range("I106").Value = 991.38464
$ws.Range("K106").Value = 991.38464
$ws.Range("M106").Value = -360.38464
$ws.Range("H107").Value = 736.1905
$ws.Range("I107").Value = 752.85
$ws.Range("J107").Value = 403
$ws.Range("K107").Value = 752.85
$ws.Range("L107").Value = 403
$ws.Range("M107").Value = 1167.15
$ws.Range("N107").Value = -4243
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882
$ws.Range("H112").Value = 46365.043
$ws.Range("J112").Value = 48455.5
$ws.Range("L112").Value = 145366.5
$ws.Range("N112").Value = -147582.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 25000
$ws.Range("I19").Value = 25000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -24771
$ws.Range("N19").ClearContents()
$ws.Range("H45").Value = 4830.5776
$ws.Range("I45").Value = 5936
$ws.Range("K45").Value = 5936
$ws.Range("M45").Value = -5559
$ws.Range("H61").Value = 4198.643
$ws.Range("I61").Value = 3552.9
$ws.Range("K61").Value = 3552.9
$ws.Range("M61").Value = -3340.9
$ws.Range("H69").Value = 500459
$ws.Range("J69").Value = 500459
$ws.Range("L69").Value = 500459
$ws.Range("N69").Value = -501957
$ws.Range("H72").Value = 500459
$ws.Range("J72").Value = 500459
$ws.Range("L72").Value = 1501377
$ws.Range("N72").Value = -1508865
$ws.Range("H136").Value = 4198.643
$ws.Range("I136").Value = 3552.9
$ws.Range("K136").Value = 10658.7
$ws.Range("M136").Value = -8108.700000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3499.889
$ws.Range("I86").Value = 3265.5
$ws.Range("J86").Value = 3968.6667
$ws.Range("K86").Value = 3265.5
$ws.Range("L86").Value = 3968.6667
$ws.Range("M86").Value = -2142.5
$ws.Range("N86").Value = -6214.6667
$ws.Range("H89").Value = 3499.889
$ws.Range("I89").Value = 3265.5
$ws.Range("J89").Value = 3968.6667
$ws.Range("K89").Value = 16327.5
$ws.Range("L89").Value = 19843.3335
$ws.Range("M89").Value = -10711.5
$ws.Range("N89").Value = -31075.3335
$ws.Range("H105").Value = 3131.375
$ws.Range("I105").Value = 3291.8462
$ws.Range("J105").Value = 2436
$ws.Range("K105").Value = 3291.8462
$ws.Range("L105").Value = 2436
$ws.Range("M105").Value = -1544.8462
$ws.Range("N105").Value = -5930
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 175000
$ws.Range("J20").Value = 175000
$ws.Range("L20").Value = 175000
$ws.Range("N20").Value = -175472
$ws.Range("H30").Value = 175000
$ws.Range("J30").Value = 175000
$ws.Range("L30").Value = 175000
$ws.Range("N30").Value = -175182
$ws.Range("H31").Value = 45423.875
$ws.Range("I31").Value = 54084.684
$ws.Range("J31").Value = 12512.8
$ws.Range("K31").Value = 54084.684
$ws.Range("L31").Value = 12512.8
$ws.Range("M31").Value = -53789.684
$ws.Range("N31").Value = -13102.8
$ws.Range("H34").Value = 45423.875
$ws.Range("I34").Value = 54084.684
$ws.Range("J34").Value = 12512.8
$ws.Range("K34").Value = 54084.684
$ws.Range("L34").Value = 12512.8
$ws.Range("M34").Value = -53882.684
$ws.Range("N34").Value = -12916.8
$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 15000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 75000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -81240
$ws.Range("H86").Value = 6699.75
$ws.Range("I86").Value = 7000
$ws.Range("J86").Value = 6599.6665
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 6599.6665
$ws.Range("M86").Value = -5877
$ws.Range("N86").Value = -8845.666499999999
$ws.Range("H89").Value = 6699.75
$ws.Range("I89").Value = 7000
$ws.Range("J89").Value = 6599.6665
$ws.Range("K89").Value = 35000
$ws.Range("L89").Value = 32998.3325
$ws.Range("M89").Value = -29384
$ws.Range("N89").Value = -44230.3325
$ws.Range("H128").Value = 175000
$ws.Range("J128").Value = 175000
$ws.Range("L128").Value = 175000
$ws.Range("N128").Value = -184960
$ws.Range("H132").Value = 3161.4285
$ws.Range("I132").Value = 3073.5757
$ws.Range("K132").Value = 9220.7271
$ws.Range("M132").Value = -6690.7271
$ws.Range("H134").Value = 10266.667
$ws.Range("I134").Value = 7263.1665
$ws.Range("J134").Value = 18276
$ws.Range("K134").Value = 21789.4995
$ws.Range("L134").Value = 54828
$ws.Range("M134").Value = -19254.4995
$ws.Range("N134").Value = -59898
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140
$ws.Range("H137").Value = 92984.28999999999
$ws.Range("J137").Value = 92984.28999999999
$ws.Range("L137").Value = 92984.28999999999
$ws.Range("N137").Value = -103184.29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1327
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 27781596
$ws.Range("I18").Value = 27781596
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 27781596
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -27781303
$ws.Range("N18").ClearContents()
$ws.Range("H102").Value = 2162.2964
$ws.Range("I102").Value = 2208.6538
$ws.Range("J102").Value = 957
$ws.Range("K102").Value = 2208.6538
$ws.Range("L102").Value = 957
$ws.Range("M102").Value = -586.6538
$ws.Range("N102").Value = -4201
$ws.Range("H122").Value = 1026.8948
$ws.Range("I122").Value = 974.4286
$ws.Range("K122").Value = 2923.2858
$ws.Range("M122").Value = -473.2857999999997
$ws.Range("H123").Value = 39900
$ws.Range("J123").Value = 39900
$ws.Range("L123").Value = 39900
$ws.Range("N123").Value = -44800
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11298.7
$ws.Range("I7").Value = 12872.75
$ws.Range("K7").Value = 12872.75
$ws.Range("M7").Value = -12760.75
$ws.Range("H16").Value = 6330.0645
$ws.Range("I16").Value = 2811.1428
$ws.Range("J16").Value = 7356.4165
$ws.Range("K16").Value = 2811.1428
$ws.Range("L16").Value = 7356.4165
$ws.Range("M16").Value = -2641.1428
$ws.Range("N16").Value = -7696.4165
$ws.Range("H40").Value = 5616.5
$ws.Range("I40").Value = 4939.3
$ws.Range("J40").Value = 9002.5
$ws.Range("K40").Value = 4939.3
$ws.Range("L40").Value = 9002.5
$ws.Range("M40").Value = -4803.3
$ws.Range("N40").Value = -9274.5
$ws.Range("H93").Value = 2498.5264
$ws.Range("I93").Value = 2445.4119
$ws.Range("J93").Value = 2950
$ws.Range("K93").Value = 2445.4119
$ws.Range("L93").Value = 2950
$ws.Range("M93").Value = -1197.4119
$ws.Range("N93").Value = -5446
$ws.Range("H126").Value = 11298.7
$ws.Range("I126").Value = 12872.75
$ws.Range("K126").Value = 38618.25
$ws.Range("M126").Value = -36148.25
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 178250
$ws.Range("J110").Value = 178250
$ws.Range("L110").Value = 178250
$ws.Range("N110").Value = -186430
$ws.Range("H126").Value = 2318.818
$ws.Range("I126").Value = 1850.7
$ws.Range("K126").Value = 5552.1
$ws.Range("M126").Value = -3082.1
